# Revert Jason's 9/13/2023 overwrite of the Encapsulation "Part2" deck:
#   - remove the 3 slides he added that day (sldId 441, 439, 440 -> the
#     last three slides in the deck: "HWVaporSalesManager",
#     "Importing HWVaporSalesManager", and the import-screenshots slide)
#   - restore the cached date-field text on every master/layout back to
#     9/12/2023 (the day the deck was last saved before his edits)

$p = $ppt.ActivePresentation

# --- 1. Drop the three trailing slides Jason added on 9/13 -----------------
$count = $p.Slides.Count
for ($i = $count; $i -ge ($count - 2); $i--) {
    $p.Slides.Item($i).Delete()
}

# --- 2. Roll the cached "datetimeFigureOut" field text back a day ----------
function Set-DatePlaceholderText($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = "9/12/2023"
        }
    }
}

# Slide master
Set-DatePlaceholderText $p.SlideMaster.Shapes

# Every slide layout under the master
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Set-DatePlaceholderText $layouts.Item($li).Shapes
}

# Notes master
Set-DatePlaceholderText $p.NotesMaster.Shapes

# Handout master
Set-DatePlaceholderText $p.HandoutMaster.Shapes
